$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("_input")

$srcRange = $ws.Range("A5:A28")
$dstRange = $ws.Range("E5:E28")
$srcRange.Copy($dstRange)

for ($r = 5; $r -le 28; $r++) {
    $hour = $r - 4
    $ws.Cells.Item($r, 5).Value = "$hour`:00"
}
